$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits (rows 2-25), before any row deletions ---

# F2: 18.03 -> blank
$ws.Range("F2").ClearContents()

# E6: blank -> -5.7
$ws.Range("E6").Value = -5.7

# E8: -6.6 -> blank
$ws.Range("E8").ClearContents()

# E18: blank -> -8.5
$ws.Range("E18").Value = -8.5

# E20: -7.2 -> blank
$ws.Range("E20").ClearContents()

# E23: blank -> -7
$ws.Range("E23").Value = -7

# E25: -7.1 -> blank
$ws.Range("E25").ClearContents()

# --- Remove rows "RM 232" (row 26) and "SC 92" (originally row 28) ---

# Delete row 26 ("RM 232"); subsequent rows shift up by one
$ws.Rows(26).Delete()

# "SC 92" was originally row 28, now sits at row 27 after the shift
$ws.Rows(27).Delete()

# --- Cell value edits for remaining rows, using the final (post-deletion) row numbers ---

# SC 101 is now row 27: C blank -> 10
$ws.Range("C27").Value = 10

# SC 105 is now row 28: C 11.1 -> blank
$ws.Range("C28").ClearContents()

# SC 119 is now row 29: C 11.2 -> blank
$ws.Range("C29").ClearContents()

# SC 120 is now row 30: C blank -> 11.4, E blank -> -5.7, F blank -> 16.89
$ws.Range("C30").Value = 11.4
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

# SC 193 is now row 32: C 10.5 -> blank
$ws.Range("C32").ClearContents()
